$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated expression/specificity values (new TPM-based numbers) for the
# Angpt1 -> Itga5 ligand-receptor pair table, columns E:T across rows 2-16.
# Row order (Sending cluster / Target cluster) is unchanged; only the
# numeric metrics were refreshed with the new TPM-derived computation.
$data = @(
    @(2,0.6666666666666666,0.1497083333333333,0.449125,0.006513369349540601,0.006769619242096868,3,1,28.85518433333334,86.56555300000001,0.1999651185353207,0.2044513327926365,4.319861554569445,38.878753991125,0.001302446674045211,0.001384057676545382),
    @(2,0.6666666666666666,0.1497083333333333,0.449125,0.006513369349540601,0.006769619242096868,3,1,51.17424933333334,153.522748,0.3546352265743414,0.3625914622481308,7.661211577277778,68.9509041955,0.002309870215036702,0.002454606139854986),
    @(2,0.6666666666666666,0.1497083333333333,0.449125,0.006513369349540601,0.006769619242096868,3,1,29.393479,88.180437,0.2036954761578358,0.2082653809291453,4.400448751958333,39.604038767625,0.001326743871046526,0.001409877330200576),
    @(2,0.6666666666666666,0.1497083333333333,0.449125,0.006513369349540601,0.006769619242096868,2,1,9.499066500000001,18.998133,0.0658280999596015,0.04486996822421697,1.4220894139375,8.532536483625002,0.0004287627286153633,0.0003037526002829342),
    @(2,0.6666666666666666,0.1497083333333333,0.449125,0.006513369349540601,0.006769619242096868,3,1,25.37910966666666,76.13732899999999,0.1758760787729007,0.1798218558058706,3.799464209680555,34.195177887125,0.001145545860796799,0.00121732549521299),
    @(3,1,20.22494433333334,60.67483300000001,0.8799278542737426,0.9145460996110527,3,1,28.85518433333334,86.56555300000001,0.1999651185353207,0.2044513327926365,583.5944968697389,5252.35047182765,0.1759548776823793,0.186980168965787),
    @(3,1,20.22494433333334,60.67483300000001,0.8799278542737426,0.9145460996110527,3,1,51.17424933333334,153.522748,0.3546352265743414,0.3625914622481308,1034.996344066787,9314.967096601085,0.3120534139694428,0.3316066075512963),
    @(3,1,20.22494433333334,60.67483300000001,0.8799278542737426,0.9145460996110527,3,1,29.393479,88.180437,0.2036954761578358,0.2082653809291453,594.4814765380024,5350.333288842022,0.1792373232608328,0.1904682918127599),
    @(3,1,20.22494433333334,60.67483300000001,0.8799278542737426,0.9145460996110527,2,1,9.499066500000001,18.998133,0.0658280999596015,0.04486996822421697,192.1180911811315,1152.708547086789,0.05792397874836959,0.0410356544291295),
    @(3,1,20.22494433333334,60.67483300000001,0.8799278542737426,0.9145460996110527,3,1,25.37910966666666,76.13732899999999,0.1758760787729007,0.1798218558058706,513.2910802378952,4619.619722141058,0.1547582606127182,0.1644553768520801),
    @(2,1,2.610123,5.220245999999999,0.1135587763767167,0.07868428114685043,3,1,28.85518433333334,86.56555300000001,0.1999651185353207,0.2044513327926365,75.31558029767301,451.893481786038,0.02270779417889612,0.01608710615030409),
    @(2,1,2.610123,5.220245999999999,0.1135587763767167,0.07868428114685043,3,1,51.17424933333334,153.522748,0.3546352265743414,0.3625914622481308,133.571085192668,801.4265111560079,0.04027194238986189,0.02853024855697953),
    @(2,1,2.610123,5.220245999999999,0.1135587763767167,0.07868428114685043,3,1,29.393479,88.180437,0.2036954761578358,0.2082653809291453,76.72059558791699,460.3235735275019,0.02313140902595651,0.01638721178618477),
    @(2,1,2.610123,5.220245999999999,0.1135587763767167,0.07868428114685043,2,1,9.499066500000001,18.998133,0.0658280999596015,0.04486996822421697,24.7937319501795,99.17492780071801,0.007475358482616539,0.003530561194804533),
    @(2,1,2.610123,5.220245999999999,0.1135587763767167,0.07868428114685043,3,1,25.37910966666666,76.13732899999999,0.1758760787729007,0.1798218558058706,66.24259786048898,397.4555871629339,0.00736268535107,0.01414915345857752)
)

$rowCount = $data.Count
$colCount = $data[0].Count
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $rowVals[$j]
    }
}

$ws.Range("E2:T16").Value = $arr
